$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format applied first, otherwise Excel would auto-convert the
# assigned string into a numeric value (and potentially drop significant
# trailing zeros), which would not match the source data (plain text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "26.077.17"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "1.651.25"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "217.07"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "0.5271"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "0.06309"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "20.32"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("D11").Value = "0.07783"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "4.517"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "1.657.44"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "1.878.28"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "0.5468"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "0.0₅8183"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "65.31"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "26.064.56"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "4.573"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "190.30"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "6.011"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "143.40"
$ws.Range("E25").Value = "  +3.26%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "7.215"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "15.97"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").Value = "1.445"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").Value = "0.05809"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "3.540"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "1.589"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "2.793"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "0.9417"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").Value = "0.5744"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "0.8486"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "104.43"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "5.713"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("D44").Value = "1.028.62"
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("D45").Value = "1.793.08"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "1.002"
$ws.Range("D48").Value = "0.4326"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "0.05138"
$ws.Range("E51").Value = "  -1.42%  "
